# svmvIncidentes_1.xlsx -- add 2022 column + running "total" column
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Header row: insert the "2022" year header in G1, and push the
#    old "total" header out to the new H1 column.
# ---------------------------------------------------------------
$ws.Range("G1").Value = 2022
$ws.Range("H1").Value = "total"

# ---------------------------------------------------------------
# 2. Replace the old running-total column (G2:G10, which held the
#    sum across 2017-2021) with the freshly reported 2022 figures.
# ---------------------------------------------------------------
$ws.Range("G2").Value = 28
$ws.Range("G3").Value = 25
$ws.Range("G4").Value = 21
$ws.Range("G5").Value = 7
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 84

# Corrected 2017 figure for the "Homicidio(s) seguido de suicidio(s)" row
$ws.Range("B6").Value = 2

# ---------------------------------------------------------------
# 3. New column H: row total across the year columns (B:G).
# ---------------------------------------------------------------
$ws.Range("H2").Formula = "=SUM(B2:G2)"
$ws.Range("H3").Formula = "=SUM(B3:G3)"
$ws.Range("H4").Formula = "=SUM(B4:G4)"
$ws.Range("H5").Formula = "=SUM(B5:G5)"
$ws.Range("H6").Formula = "=SUM(B6:G6)"
$ws.Range("H7").Formula = "=SUM(B7:G7)"
$ws.Range("H8").Formula = "=SUM(B8:G8)"
$ws.Range("H9").Formula = "=SUM(B9:G9)"
$ws.Range("H10").Formula = "=SUM(B10:G10)"

# ---------------------------------------------------------------
# 4. Row 10 ("Total de victimas mujeres"): turn the hard-coded
#    numbers into column sums (rows 2-8), matching the rest of the
#    sheet.
# ---------------------------------------------------------------
$ws.Range("B10").Formula = "=SUM(B2:B8)"
$ws.Range("C10").Formula = "=SUM(C2:C8)"
$ws.Range("D10").Formula = "=SUM(D2:D8)"
$ws.Range("E10").Formula = "=SUM(E2:E8)"
$ws.Range("F10").Formula = "=SUM(F2:F8)"
$ws.Range("G10").Formula = "=SUM(G2:G8)"

# ---------------------------------------------------------------
# 5. Wording tweak: shorten "más de dos víctimas" to "+ dos víctimas".
# ---------------------------------------------------------------
$ws.Range("A6").Value = "Homicidio(s) seguido de suicidio(s), (+ dos víctimas)"

# ---------------------------------------------------------------
# 6. Formatting: center (horizontal + vertical) every data cell in
#    the new B1:H10 block.
# ---------------------------------------------------------------
$data = $ws.Range("B1:H10")
$data.HorizontalAlignment = -4108
$data.VerticalAlignment = -4108

# Keep the "Aptos Narrow" font on the totals row (F9:G9 already used
# it before this edit).
$ws.Range("F9:G9").Font.Name = "Aptos Narrow"

# ---------------------------------------------------------------
# 7. Selection cosmetics - matches where the author's cursor ended
#    up after entering the new row of data.
# ---------------------------------------------------------------
$ws.Range("A11").Select()
